# Generate Report for Handback
# Updates the "3b55e8f2-8a12-4a9f-89d3-9a04209009b2" row (row 7) on both the
# zh-cn and de-de language sheets: the handback for that file is now in, but
# it is not the newest version of the handback, so a new "Latest Target
# File" / "Latest Handback File" / "Latest Handback DateTime" / Error Detail
# are recorded, together with a hyperlink on the new "Latest Target File"
# cell (column I).

$wb = $excel.ActiveWorkbook

$handbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2ea6a193ab8f3da1538c3c755f7722e7cad77ba1/e2e/3b55e8f2-8a12-4a9f-89d3-9a04209009b2.md"
$handbackDisplay = "3b55e8f2-8a12-4a9f-89d3-9a04209009b2.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bf7f2ffcc290d43d5322a26e09dbcf3c49d94ce3/e2e/3b55e8f2-8a12-4a9f-89d3-9a04209009b2.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2ea6a193ab8f3da1538c3c755f7722e7cad77ba1/e2e/3b55e8f2-8a12-4a9f-89d3-9a04209009b2.md."
$hyperlinkColor = 15570276

# ---------------------------------------------------------------------
# zh-cn sheet, row 7
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("J7").Value2 = "3b55e8f2-8a12-4a9f-89d3-9a04209009b2.e2df93ca23dcfa0d807ae5edf2ae828caf71a005.zh-cn.xlf"
$wsZh.Range("K7").Value2 = "2016-08-24 22:55:46"
$wsZh.Range("P7").Value2 = $errorDetail

$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $handbackUrl, "", "", $handbackDisplay) | Out-Null
$wsZh.Range("I7").Font.Color = $hyperlinkColor

# ---------------------------------------------------------------------
# de-de sheet, row 7
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("J7").Value2 = "3b55e8f2-8a12-4a9f-89d3-9a04209009b2.e2df93ca23dcfa0d807ae5edf2ae828caf71a005.de-de.xlf"
$wsDe.Range("K7").Value2 = "2016-08-24 22:55:53"
$wsDe.Range("P7").Value2 = $errorDetail

$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $handbackUrl, "", "", $handbackDisplay) | Out-Null
$wsDe.Range("I7").Font.Color = $hyperlinkColor
